# Split the run containing "su número de casilla" into two runs:
#   "su"  and  " número de casilla"
# (matches the target OOXML diff, which breaks one <w:r> into two).

$d = $word.ActiveDocument

$needle = "su número de casilla"
$splitAfter = "su"

foreach ($p in $d.Paragraphs) {
    $full = $p.Range.Text
    if ($full -ne $null -and $full.Length -gt 0) {
        $idx = $full.IndexOf($needle)
        if ($idx -ge 0) {
            $paraStart = $p.Range.Start
            $splitPoint = $paraStart + $idx + $splitAfter.Length

            # Toggling a character-formatting property on the sub-range forces
            # the engine to break the run in two at that boundary while
            # leaving the visible formatting unchanged (Bold is immediately
            # turned back off).
            $firstPart = $d.Range($paraStart + $idx, $splitPoint)
            $firstPart.Bold = 1
            $firstPart.Bold = 0

            break
        }
    }
}
